$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4549446666666667
$ws.Range("H2").Value = 1.364834
$ws.Range("I2").Value = 0.8656500014587819
$ws.Range("J2").Value = 0.8656500014587818
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2773296666666667
$ws.Range("N2").Value = 0.8319890000000001
$ws.Range("Q2").Value = 0.1261696527584445
$ws.Range("R2").Value = 1.135526874826
$ws.Range("S2").Value = 0.8656500014587819
$ws.Range("T2").Value = 0.8656500014587818

# Row 3
$ws.Range("I3").Value = 0.02407814503842938
$ws.Range("J3").Value = 0.02407814503842938
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2773296666666667
$ws.Range("N3").Value = 0.8319890000000001
$ws.Range("Q3").Value = 0.003509422045222222
$ws.Range("R3").Value = 0.031584798407
$ws.Range("S3").Value = 0.02407814503842938
$ws.Range("T3").Value = 0.02407814503842938

# Row 4
$ws.Range("G4").Value = 0.05795366666666666
$ws.Range("H4").Value = 0.173861
$ws.Range("I4").Value = 0.1102718535027888
$ws.Range("J4").Value = 0.1102718535027888
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2773296666666667
$ws.Range("N4").Value = 0.8319890000000001
$ws.Range("Q4").Value = 0.01607227105877778
$ws.Range("R4").Value = 0.144650439529
$ws.Range("S4").Value = 0.1102718535027888
$ws.Range("T4").Value = 0.1102718535027888
